# This workbook tracks weekly "Ajo" (garlic) price observations at
# Vega Monumental Concepción. A new weekly observation is inserted as
# row 19 (pushing the existing rows 19-87 down to rows 20-88), and the
# new row is populated with this week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 19; this shifts rows 19..87 down to 20..88
# and carries the existing formatting/styles (e.g. the date number
# format on column D) down with them.
$ws.Rows.Item(19).Insert()

# Populate the newly inserted row 19 with the new weekly observation.
$ws.Cells.Item(19, 1).Value = 11
$ws.Cells.Item(19, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(19, 3).Value = "Bíobío"
$ws.Cells.Item(19, 4).Value = "09/14/2021"
$ws.Cells.Item(19, 5).Value = 8
$ws.Cells.Item(19, 6).Value = 100112003
$ws.Cells.Item(19, 7).Value = "Ajo"
$ws.Cells.Item(19, 8).Value = "Chino"
$ws.Cells.Item(19, 9).Value = "Primera"
$ws.Cells.Item(19, 10).Value = 400
$ws.Cells.Item(19, 11).Value = 16000
$ws.Cells.Item(19, 12).Value = 17000
$ws.Cells.Item(19, 13).Value = 16500
$ws.Cells.Item(19, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(19, 15).Value = "China"
$ws.Cells.Item(19, 16).Value = 1650
$ws.Cells.Item(19, 17).Value = 10
$ws.Cells.Item(19, 18).Value = "Hortaliza"
